$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Inscritos" (E), "Pagos" (F) and "Inscrições homologadas" (H) counts
$ws.Range("E2").Value = 84
$ws.Range("F2").Value = 58
$ws.Range("H2").Value = 58

$ws.Range("E5").Value = 112
$ws.Range("F5").Value = 70
$ws.Range("H5").Value = 70

$ws.Range("E6").Value = 38

$ws.Range("E10").Value = 388
$ws.Range("F10").Value = 175
$ws.Range("H10").Value = 175

$ws.Range("E11").Value = 258

$ws.Range("E12").Value = 369
$ws.Range("F12").Value = 204
$ws.Range("H12").Value = 204

$ws.Range("E14").Value = 97
$ws.Range("F14").Value = 47
$ws.Range("H14").Value = 47

$ws.Range("E15").Value = 123

$ws.Range("E16").Value = 164

$ws.Range("E17").Value = 71

$ws.Range("E18").Value = 47

$ws.Range("E21").Value = 121

$ws.Range("E22").Value = 143

$ws.Range("E23").Value = 166

$ws.Range("E24").Value = 165

$ws.Range("E25").Value = 193

$ws.Range("E27").Value = 256

$ws.Range("E28").Value = 151

$ws.Range("E30").Value = 166
$ws.Range("F30").Value = 92
$ws.Range("H30").Value = 92

$ws.Range("E32").Value = 154
$ws.Range("F32").Value = 83
$ws.Range("H32").Value = 83

$ws.Range("E33").Value = 235

$ws.Range("E34").Value = 170
$ws.Range("F34").Value = 100
$ws.Range("H34").Value = 100

$ws.Range("E35").Value = 112

$ws.Range("E36").Value = 49

$ws.Range("E37").Value = 124
$ws.Range("F37").Value = 60
$ws.Range("H37").Value = 60

$ws.Range("E38").Value = 78

$ws.Range("E39").Value = 156

$ws.Range("E40").Value = 213
$ws.Range("F40").Value = 89
$ws.Range("H40").Value = 89

$ws.Range("E41").Value = 313
$ws.Range("F41").Value = 131
$ws.Range("H41").Value = 131

$ws.Range("E42").Value = 279
$ws.Range("F42").Value = 146
$ws.Range("H42").Value = 146

$ws.Range("E43").Value = 94

$ws.Range("E44").Value = 250

$ws.Range("E46").Value = 244
$ws.Range("F46").Value = 131
$ws.Range("H46").Value = 131

$ws.Range("E47").Value = 352
$ws.Range("F47").Value = 165
$ws.Range("H47").Value = 165

$ws.Range("E48").Value = 163
$ws.Range("F48").Value = 64
$ws.Range("H48").Value = 64

$ws.Range("E49").Value = 233

$ws.Range("E50").Value = 206
$ws.Range("F50").Value = 73
$ws.Range("H50").Value = 73

$ws.Range("E51").Value = 190
$ws.Range("F51").Value = 76
$ws.Range("H51").Value = 76
